$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (SWLS_DRS_Composite_Z, .GOATTotalScore_Z)
$ws.Range("C8").Value = 0.14814730973738
$ws.Range("D8").Value = 0.01810307127738879
$ws.Range("E8").Value = 0.1240665813995394
$ws.Range("F8").Value = 0.1736406001703299
$ws.Range("G8").Value = 0.134271690424067
$ws.Range("H8").Value = 0.1596206999675694

# Row 9 (SWLS_DRS_Composite_Z, PTA_Reverse_Z)
$ws.Range("C9").Value = 0.2606213405587098
$ws.Range("D9").Value = 0.01001885071711333
$ws.Range("E9").Value = 0.232093937381568
$ws.Range("F9").Value = 0.2747122803586632
$ws.Range("G9").Value = 0.2515655372898408
$ws.Range("H9").Value = 0.2665331144783281

# Row 10 (SWLS_DRS_Composite_Z, WorsT.GCS.and.Pupils.GCSTotalScore_Z)
$ws.Range("C10").Value = 0.242901726420291
$ws.Range("D10").Value = 0.01520806452611467
$ws.Range("E10").Value = 0.2022570932042213
$ws.Range("F10").Value = 0.2725805542232504
$ws.Range("G10").Value = 0.2339862158571543
$ws.Range("H10").Value = 0.2535812536075055

$wb.Save()
